$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 68.42856999999999
$ws.Range("I5").Value = 68.42856999999999
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 68.42856999999999
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 46.57143000000001
$ws.Range("N5").ClearContents()
# Row 15
$ws.Range("H15").Value = 2407.2046
$ws.Range("I15").Value = 2407.2046
$ws.Range("K15").Value = 7221.6138
$ws.Range("M15").Value = -7052.6138
# Row 17
$ws.Range("H17").Value = 1044770.94
$ws.Range("J17").Value = 1305844.9
$ws.Range("L17").Value = 3917534.7
$ws.Range("N17").Value = -3917870.7
# Row 21
$ws.Range("H21").Value = 19142.857
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19532
# Row 23
$ws.Range("H23").Value = 19142.857
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19766
# Row 62
$ws.Range("H62").Value = 6420.35
$ws.Range("I62").Value = 4066.7778
$ws.Range("K62").Value = 4066.7778
$ws.Range("M62").Value = -3442.7778
# Row 65
$ws.Range("H65").Value = 6420.35
$ws.Range("I65").Value = 4066.7778
$ws.Range("K65").Value = 20333.889
$ws.Range("M65").Value = -17213.889
# Row 80
$ws.Range("H80").Value = 1591.95
$ws.Range("I80").Value = 887.375
$ws.Range("K80").Value = 2662.125
$ws.Range("M80").Value = -1664.125
# Row 83
$ws.Range("H83").Value = 1591.95
$ws.Range("I83").Value = 887.375
$ws.Range("K83").Value = 7986.375
$ws.Range("M83").Value = -2994.375
# Row 86
$ws.Range("H86").Value = 4981.6665
$ws.Range("J86").Value = 5999.8335
$ws.Range("L86").Value = 5999.8335
$ws.Range("N86").Value = -8245.833500000001
# Row 89
$ws.Range("H89").Value = 4981.6665
$ws.Range("J89").Value = 5999.8335
$ws.Range("L89").Value = 29999.1675
$ws.Range("N89").Value = -41231.1675
# Row 116
$ws.Range("H116").Value = 8069.2
$ws.Range("I116").Value = 7438.6
$ws.Range("K116").Value = 7438.6
$ws.Range("M116").Value = -3996.6
# Row 137
$ws.Range("H137").Value = 3320.5217
$ws.Range("I137").Value = 1844.3636
$ws.Range("K137").Value = 5533.0908
$ws.Range("M137").Value = -2983.0908

$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 47368772
$ws.Range("I17").Value = 355.8125
$ws.Range("J17").Value = 300000320
$ws.Range("K17").Value = 355.8125
$ws.Range("L17").Value = 300000320
$ws.Range("M17").Value = -182.8125
$ws.Range("N17").Value = -300000666
# Row 135
$ws.Range("H135").Value = 62499
$ws.Range("J135").Value = 62499
$ws.Range("L135").Value = 62499
$ws.Range("N135").Value = -72639
# Row 139
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1188.375
$ws.Range("I94").Value = 1219.619
$ws.Range("J94").Value = 969.6667
$ws.Range("K94").Value = 1219.619
$ws.Range("L94").Value = 969.6667
$ws.Range("M94").Value = -768.6189999999999
$ws.Range("N94").Value = -1871.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 714470.9
$ws.Range("J4").Value = 714470.9
$ws.Range("L4").Value = 714470.9
$ws.Range("N4").Value = -714694.9
# Row 7
$ws.Range("H7").Value = 325.33334
$ws.Range("I7").Value = 157.25
$ws.Range("K7").Value = 157.25
$ws.Range("M7").Value = -44.25
# Row 11
$ws.Range("H11").Value = 140.07143
$ws.Range("I11").Value = 123.72727
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 123.72727
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = 16.27273
$ws.Range("N11").Value = -480
# Row 31
$ws.Range("H31").Value = 32766.756
$ws.Range("I31").Value = 3206.842
$ws.Range("K31").Value = 3206.842
$ws.Range("M31").Value = -2911.842
# Row 34
$ws.Range("H34").Value = 32766.756
$ws.Range("I34").Value = 3206.842
$ws.Range("K34").Value = 3206.842
$ws.Range("M34").Value = -3004.842
# Row 62
$ws.Range("H62").Value = 11053.909
$ws.Range("J62").Value = 14350
$ws.Range("L62").Value = 14350
$ws.Range("N62").Value = -15598
# Row 65
$ws.Range("H65").Value = 11053.909
$ws.Range("J65").Value = 14350
$ws.Range("L65").Value = 71750
$ws.Range("N65").Value = -77990
# Row 105
$ws.Range("H105").Value = 10595.6
$ws.Range("I105").Value = 8491.75
$ws.Range("K105").Value = 8491.75
$ws.Range("M105").Value = -6744.75
# Row 132
$ws.Range("H132").Value = 3515.611
$ws.Range("I132").Value = 3491.4375
$ws.Range("K132").Value = 10474.3125
$ws.Range("M132").Value = -7944.3125

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 3154806.2
$ws.Range("I4").Value = 5606471.5
$ws.Range("J4").Value = 703141.4
$ws.Range("K4").Value = 16819414.5
$ws.Range("L4").Value = 2109424.2
$ws.Range("M4").Value = -16819302.5
$ws.Range("N4").Value = -2109648.2
# Row 57
$ws.Range("H57").Value = 3999.5
$ws.Range("I57").Value = 3999.5
$ws.Range("K57").Value = 11998.5
$ws.Range("M57").Value = -11439.5
# Row 64
$ws.Range("H64").Value = 100006830
$ws.Range("J64").Value = 7836.3335
$ws.Range("L64").Value = 23509.0005
$ws.Range("N64").Value = -24049.0005
# Row 67
$ws.Range("H67").Value = 100006830
$ws.Range("J67").Value = 7836.3335
$ws.Range("L67").Value = 23509.0005
$ws.Range("N67").Value = -25381.0005
# Row 129
$ws.Range("H129").Value = 3789682.5
$ws.Range("I129").Value = 527.7273
$ws.Range("J129").Value = 7578837.5
$ws.Range("K129").Value = 1583.1819
$ws.Range("L129").Value = 22736512.5
$ws.Range("M129").Value = 3416.8181
$ws.Range("N129").Value = -22746512.5

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 113
$ws.Range("H113").Value = 5254.55
$ws.Range("I113").Value = 2051.5
$ws.Range("J113").Value = 8457.6
$ws.Range("K113").Value = 2051.5
$ws.Range("L113").Value = 8457.6
$ws.Range("M113").Value = 118.5
$ws.Range("N113").Value = -12797.6
# Row 126
$ws.Range("H126").Value = 5406.263
$ws.Range("J126").Value = 9419.625
$ws.Range("L126").Value = 28258.875
$ws.Range("N126").Value = -33198.875
# Row 132
$ws.Range("H132").Value = 11839.421
$ws.Range("I132").Value = 3492.5
$ws.Range("J132").Value = 17909.908
$ws.Range("K132").Value = 10477.5
$ws.Range("L132").Value = 53729.724
$ws.Range("M132").Value = -7947.5
$ws.Range("N132").Value = -58789.724
# Row 133
$ws.Range("H133").Value = 69993.5
$ws.Range("J133").Value = 69993.5
$ws.Range("L133").Value = 69993.5
$ws.Range("N133").Value = -80113.5

$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 1671.7142
$ws.Range("I19").Value = 283.83334
$ws.Range("K19").Value = 283.83334
$ws.Range("M19").Value = -113.83334
# Row 46
$ws.Range("H46").Value = 2753.3572
$ws.Range("J46").Value = 3012.4167
$ws.Range("L46").Value = 3012.4167
$ws.Range("N46").Value = -3388.4167
# Row 93
$ws.Range("H93").Value = 22540.344
$ws.Range("I93").Value = 17972.38
$ws.Range("J93").Value = 29392.285
$ws.Range("K93").Value = 17972.38
$ws.Range("L93").Value = 29392.285
$ws.Range("M93").Value = -16724.38
$ws.Range("N93").Value = -31888.285
# Row 122
$ws.Range("H122").Value = 5225.579
$ws.Range("I122").Value = 4296.5
$ws.Range("K122").Value = 12889.5
$ws.Range("M122").Value = -10439.5

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 65000
$ws.Range("J46").Value = 65000
$ws.Range("L46").Value = 65000
$ws.Range("N46").Value = -65462
# Row 81
$ws.Range("H81").Value = 1495
$ws.Range("I81").Value = 1495
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2990
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1929
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 1495
$ws.Range("I84").Value = 1495
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 14950
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -9646
$ws.Range("N84").ClearContents()
# Row 100
$ws.Range("H100").Value = 788.2353000000001
$ws.Range("I100").Value = 785.9286
$ws.Range("K100").Value = 1571.8572
$ws.Range("M100").Value = -1030.8572
# Row 107
$ws.Range("H107").Value = 1042.1428
$ws.Range("I107").Value = 1099.6666
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 3298.9998
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -1378.9998
$ws.Range("N107").Value = -6837
# Row 113
$ws.Range("H113").Value = 1268.3928
$ws.Range("I113").Value = 1322
$ws.Range("J113").Value = 1071.8334
$ws.Range("K113").Value = 3966
$ws.Range("L113").Value = 3215.5002
$ws.Range("M113").Value = -1796
$ws.Range("N113").Value = -7555.5002
# Row 126
$ws.Range("H126").Value = 5041.0386
$ws.Range("I126").Value = 4215
$ws.Range("K126").Value = 12645
$ws.Range("M126").Value = -10175
# Row 132
$ws.Range("H132").Value = 9046.25
$ws.Range("I132").Value = 6626.067
$ws.Range("K132").Value = 19878.201
$ws.Range("M132").Value = -17348.201
# Row 134
$ws.Range("H134").Value = 65000
$ws.Range("J134").Value = 65000
$ws.Range("L134").Value = 195000
$ws.Range("N134").Value = -200070
# Row 136
$ws.Range("H136").Value = 4890.95
$ws.Range("I136").Value = 2765.1333
$ws.Range("J136").Value = 11268.4
$ws.Range("K136").Value = 8295.3999
$ws.Range("L136").Value = 33805.2
$ws.Range("M136").Value = -5745.3999
$ws.Range("N136").Value = -38905.2
# Row 141
$ws.Range("H141").Value = 275715
$ws.Range("J141").Value = 275715
$ws.Range("L141").Value = 275715
$ws.Range("N141").Value = -286075
